$wb = $excel.ActiveWorkbook

# Locate the last existing sheet ("2025-07-17") and add the new sheet right after it
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-07-18"

# Copy header formatting (bold font, borders, centered alignment) from the
# previous day's sheet so the new header row matches the existing style.
$srcHeader = $wb.Worksheets.Item("2025-07-17").Range("A1:D1")
$srcHeader.Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Header row values
$newSheet.Cells.Item(1, 1).Value = "rank"
$newSheet.Cells.Item(1, 2).Value = "title"
$newSheet.Cells.Item(1, 3).Value = "author"
$newSheet.Cells.Item(1, 4).Value = "latest_episode"

# Ranking data for 2025-07-18 ("rank|title|author|latest_episode" per line)
$data = @'
1|生徒会にも穴はある！|むちまろ|第131話	ありす大ピンチ！（デジャブ編）
2|落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～|村上よしゆき 茨木野 あるてら|第４０話　勇者、聖女と元聖騎士と再会し、魚人を追っ払う（５）
3|ライドンキング|馬場康誌|第81話 大統領と失われた神器（前編）
4|ライブダンジョン！|ことりりょう(作画) dy冷凍(原作) Mika Pikazo(キャラクター原案)|第88話前半
5|最弱貴族に転生したので悪役たちを集めてみた|空野進 sorani ファルまろ|第10話　最弱貴族、部下を信じる（２）
6|35歳独身山田、異世界村に理想のセカンドハウスを作りたい　～異世界と現実のいいとこどりライフ～|出雲大吉(原作) 西尾洋一(作画) ゆのひと(キャラクター原案)|第2話①
7|デスゲームに巻き込まれた山本さん、気ままにゲームバランスを崩壊させる|ぽち(原作) カモトタツヤ(作画) 久賀フーナ(キャラクター原案)|休載イラスト
8|極振り拒否して手探りスタート！　特化しないヒーラー、仲間と別れて旅に出る|蒼井一秀(作画) 刻一(原作) MIYA*KI(キャラクター原案)|第64話
9|異世界おじさん|殆ど死んでいる(著者)|【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！
10|生徒会役員共|氏家ト全|#405
11|貧乏騎士に嫁入りしたはずが!? ~野人令嬢は皇太子妃になっても竜を狩りたい~|漫画：夏川そぞろ 原作：宮前葵 キャラクター原案：ののまろ|第10話④皇太子妃（仮）
12|【悲報】清楚系で売っていた底辺配信者、うっかり配信を切り忘れたままSS級モンスターを拳で殴り飛ばしてしまう|アトハ NEO草野 pupps|第５話　【悲報】ゆきのんとのコラボ配信、ワンパン撃破！？（１）
13|となりの席のヤツがそういう目で見てくる|mmk|第41話 教えたくない？
14|規格外のダンジョン攻略者、実は異世界帰りの元勇者|作画：やまざき君 原作：榊与一|第4話(2)
15|町人Aは悪役令嬢をどうしても救いたい　～どぶと空と氷の姫君～|原作：一色孝太郎・Parum 漫画：目黒三吉|第38話 攻撃開始
16|継母の連れ子が元カノだった|草壁レイ(作画) 紙城境介(原作) たかやKi(キャラクター原案)|第61話前半
17|十年目、帰還を諦めた転移者はいまさら主人公になる|原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう|【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！
18|異世界に転移したら山の中だった。反動で強さよりも快適さを選びました。|原作/じゃがバター(ツギクルブックス) 漫画/蔦屋空 キャラクター原案/岩崎美奈子|第32話
19|元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～|沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)|第76話その2
20|ギャルゲーマーに褒められたい|げしゅまろ(著者)|44話
21|りゅうとあまがみ|角丸柴朗(著者)|第一話・魚嫌いのウィル
22|異世界魔王と召喚少女の奴隷魔術|原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大|第126話　戦争を終わらせてみるⅡ（後編）
23|実は俺、最強でした？|原作：澄守 彩 漫画：高橋 愛|おまけ63
24|ルパン三世 異世界の姫君（ネイバーワールドプリンセス）|モンキー・パンチ／エム・ピー・ワークス 内々けやき 佐伯庸介 白狼|第101話：マッスルタの魔犬❶
25|直径3cmの召喚陣<リミットリング>で「雑魚すら呼べない」と蔑まれた底辺召喚士が頂点に立つまで|作画：まっつー 原作：空松蓮司|第4話(2)
26|やり込んだ乙女ゲームの悪役モブですが、断罪は嫌なので真っ当に生きます@COMIC|戸張ちょも（漫画） MIZUNA（原作） Ruki（キャラクター原案）|第18話
27|僕のカノジョ先生|星河蟹(作画) 孟倫（ＳＤｗｉｎｇ）(構成) 鏡遊(原作) おりょう(キャラクター原案)|75時間目-1
28|聖者無双|漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime|第90話　研究者や技術者の故郷（後半）
29|独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～|漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき|第32話 独身貴族は森で写真を撮る（1）
30|FPSゲームのコーチを引き受けたら依頼主が人気VTuberの美少女だった|漫画：真広吏希 原作：すかいふぁーむ キャラクター原案：みすみ|第8話
31|王子様の友達|すけろく(著者)|【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！
32|かつての暗殺者は来世で違う生き方をする|ツネ(漫画) 丘野優(原作) つなかわ(キャラクター原案)|第4話②
33|君の刀が折れるまで ~月宮まつりの恋難き~|イノウエ|第38話 家族のこと
34|ゲーム内最強の『裏ボス』に転生したので、主人公の代わりに最速クリアを目指します！|作画：こめぐ 原作：迅空也|第4話(2)
35|金色の文字使い ―勇者四人に巻き込まれたユニークチート―|尾崎祐介(作画) 十本スイ(原作) すまき俊悟(キャラクター原案)|第119話
36|冒険者アル あいつの魔法はおかしい@COMIC|山﨑と子（漫画） れもん（原作） sime（キャラクター原案）|第14話 ②
37|時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―|光永康則|第６６話『六花停止』③
38|魔法少女リリカルなのは EXCEEDS|都築真紀 川上修一|第４話①
39|勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが|絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)|第3話 後編
40|地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。|マツモトケンゴ|第６０話　じゃない方の戦いが始まった（２）
41|転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～|zunta(作画) はらわたさいぞう(原作)|第30話：一秒の奪い合い③
42|黒魔法寮の三悪人|斎藤キミオ|第41話 星家vs善寿賀家!だよ♪
43|オタクに優しいギャルはいない!?|のりしろちゃん 魚住さかな|【#148】氷の笑み
44|カナン様はあくまでチョロい|nonco|第142話	アミのパジャマパーティ女子会
45|ひげを剃る。そして女子高生を拾う。|足立いまる(漫画) しめさば(原作) ぶーた(キャラクター原案)|第65話その1
46|学食ガール|杉本 萌|＃01タルタルソースとコールスロー
47|無能と蔑まれた【勇気ある者】が、実は伝説の【勇者】でした~500年前の仲間＆魔王の娘と一緒に今度こそ世界を救います~|作画：シグチサトシ 原作：やすくん|第4話(2)
48|脇役に転生したはずが、いつの間にか伝説の錬金術師になってた～仲間たちが英雄でも俺は支援職なんだが～|神無月みり 相野 仁|第２６話　脇役、目論見を見破る（１）
49|すべての人類を破壊する。それらは再生できない。|横田卓馬(漫画) 伊瀬勝良(原作)|第66話その3
50|カッコウの許嫁|吉河美希|第255話	二人きりになれる場所…？
'@

$lines = $data -split "`n"
$r = 2
foreach ($line in $lines) {
    $fields = $line -split "\|", 4
    $newSheet.Cells.Item($r, 1).Value = [int]$fields[0]
    $newSheet.Cells.Item($r, 2).Value = $fields[1]
    $newSheet.Cells.Item($r, 3).Value = $fields[2]
    $newSheet.Cells.Item($r, 4).Value = $fields[3]
    $r++
}

Write-Output "Added sheet 2025-07-18 with $($r - 2) data rows"
